# Veri tabanı yedeği güncellendi
# Append a new stok (stock) row for "M5 FIBERLI SOMUN" to the SOMUNLAR sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 30): only columns A & B are populated, same as a
#     freshly-typed row (no inherited cell style from the rows above). ---
$ws.Range("A30").Value = "150.05.0506.00010"
$ws.Range("B30").Value = "M5 FIBERLI SOMUN"

# --- Column widths: best-fit widths for the data now in the sheet. The
#     target (Excel-calculated) best-fit widths include ~5/6 character of
#     padding baked into the stored width; feed ColumnWidth net of that
#     padding so the exported <col> widths land as close as possible. ---
$ws.Columns.Item(1).ColumnWidth  = 16.0221354166667   # A - Stok Kodu
$ws.Columns.Item(2).ColumnWidth  = 34.7369791666667   # B - Stok Adı
$ws.Columns.Item(3).ColumnWidth  = 7.16666666666667   # C
$ws.Columns.Item(4).ColumnWidth  = 7.16666666666667   # D
$ws.Columns.Item(5).ColumnWidth  = 4.87760416666667   # E
$ws.Columns.Item(6).ColumnWidth  = 6.87760416666667   # F
$ws.Columns.Item(7).ColumnWidth  = 7.16666666666667   # G
$ws.Columns.Item(8).ColumnWidth  = 4.87760416666667   # H
$ws.Columns.Item(9).ColumnWidth  = 6.87760416666667   # I
$ws.Columns.Item(10).ColumnWidth = 10.5924479166667   # J
$ws.Columns.Item(11).ColumnWidth = 9.59244791666667   # K

# --- View state: scroll so row 13 is at the top and select the next empty
#     row beneath the data, matching where the user left off editing. ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A31").Select()
